# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces text storage so price strings such as '1.007' or
# "14.70" keep their exact digits/trailing zeros instead of being coerced to
# numbers by the COM Value setter (matches the original inlineStr cells).

# Row 2
$ws.Range('D2').Value = '''25.178.23'
$ws.Range('E2').Value = '  -2.58%  '

# Row 3
$ws.Range('D3').Value = '''1.655.48'
$ws.Range('E3').Value = '  -4.74%  '

# Row 4
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  +0.37%  '

# Row 5
$ws.Range('D5').Value = '''236.84'
$ws.Range('E5').Value = '  -1.54%  '

# Row 6
$ws.Range('D6').Value = '''1.004'
$ws.Range('E6').Value = '  +0.39%  '

# Row 7
$ws.Range('D7').Value = '''0.4788'
$ws.Range('E7').Value = '  -8.37%  '

# Row 8
$ws.Range('D8').Value = '''0.2614'
$ws.Range('E8').Value = '  -4.90%  '

# Row 9
$ws.Range('D9').Value = '''0.05984'
$ws.Range('E9').Value = '  -2.92%  '

# Row 10
$ws.Range('D10').Value = '''0.07091'
$ws.Range('E10').Value = '  -1.44%  '

# Row 11
$ws.Range('D11').Value = '''1.674.86'
$ws.Range('E11').Value = '  -3.71%  '

# Row 12
$ws.Range('D12').Value = '''0.6200'
$ws.Range('E12').Value = '  -3.46%  '

# Row 13
$ws.Range('D13').Value = '''14.38'
$ws.Range('E13').Value = '  -3.67%  '

# Row 14
$ws.Range('D14').Value = '''4.593'
$ws.Range('E14').Value = '  -0.66%  '

# Row 15
$ws.Range('D15').Value = '''72.91'
$ws.Range('E15').Value = '  -5.92%  '

# Row 16
$ws.Range('D16').Value = '''1.006'
$ws.Range('E16').Value = '  +0.53%  '

# Row 17
$ws.Range('D17').Value = '''1.003'
$ws.Range('E17').Value = '  +0.27%  '

# Row 18
$ws.Range('D18').Value = '''25.188.40'
$ws.Range('E18').Value = '  -2.65%  '

# Row 19
$ws.Range('D19').Value = '''11.38'
$ws.Range('E19').Value = '  -2.74%  '

# Row 20
$ws.Range('D20').Value = '''0.000006550'
$ws.Range('E20').Value = '  -3.17%  '

# Row 21
$ws.Range('D21').Value = '''1.894.59'
$ws.Range('E21').Value = '  -3.74%  '

# Row 22
$ws.Range('D22').Value = '''4.421'
$ws.Range('E22').Value = '  +3.26%  '

# Row 23
$ws.Range('D23').Value = '''8.516'
$ws.Range('E23').Value = '  -1.45%  '

# Row 24
$ws.Range('D24').Value = '''5.260'
$ws.Range('E24').Value = '  -0.41%  '

# Row 25
$ws.Range('D25').Value = '''133.14'
$ws.Range('E25').Value = '  -3.78%  '

# Row 26
$ws.Range('D26').Value = '''14.70'
$ws.Range('E26').Value = '  -3.36%  '

# Row 27
$ws.Range('D27').Value = '''1.371'
$ws.Range('E27').Value = '  -9.54%  '

# Row 28
$ws.Range('D28').Value = '''1.709'
$ws.Range('E28').Value = '  -3.29%  '

# Row 29
$ws.Range('E29').Value = '  -3.27%  '

# Row 30
$ws.Range('E30').Value = '  -2.22%  '

# Row 31
$ws.Range('D31').Value = '''0.07876'
$ws.Range('E31').Value = '  -4.97%  '

# Row 32
$ws.Range('D32').Value = '''3.519'
$ws.Range('E32').Value = '  -4.38%  '

# Row 33
$ws.Range('D33').Value = '''0.04595'
$ws.Range('E33').Value = '  -0.80%  '

# Row 34
$ws.Range('D34').Value = '''2.620'
$ws.Range('E34').Value = '  -0.82%  '

# Row 35
$ws.Range('E35').Value = '  -4.89%  '

# Row 36
$ws.Range('D36').Value = '''0.5829'
$ws.Range('E36').Value = '  -5.89%  '

# Row 37
$ws.Range('D37').Value = '''2.600'
$ws.Range('E37').Value = '  -2.99%  '

# Row 38
$ws.Range('D38').Value = '''0.01543'
$ws.Range('E38').Value = '  -3.89%  '

# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''0.8440'
$ws.Range('E39').Value = '  +14.16%  '

# Row 40
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '''1.007'
$ws.Range('E40').Value = '  +0.81%  '

# Row 41
$ws.Range('D41').Value = '''1.846'
$ws.Range('E41').Value = '  -4.78%  '

# Row 42
$ws.Range('D42').Value = '''98.24'
$ws.Range('E42').Value = '  +0.27%  '

# Row 43
$ws.Range('D43').Value = '''0.3704'
$ws.Range('E43').Value = '  -3.68%  '

# Row 44
$ws.Range('D44').Value = '''4.847'
$ws.Range('E44').Value = '  -3.08%  '

# Row 45
$ws.Range('D45').Value = '''0.1133'
$ws.Range('E45').Value = '  +0.17%  '

# Row 46
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '''6.044'
$ws.Range('E46').Value = '  -3.09%  '

# Row 47
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.05151'
$ws.Range('E47').Value = '  -1.76%  '

# Row 48
$ws.Range('D48').Value = '''52.64'
$ws.Range('E48').Value = '  -3.89%  '

# Row 49
$ws.Range('D49').Value = '''29.47'
$ws.Range('E49').Value = '  -3.34%  '

# Row 50
$ws.Range('D50').Value = '''1.004'
$ws.Range('E50').Value = '  +0.34%  '

# Row 51
$ws.Range('D51').Value = '''7.355'
$ws.Range('E51').Value = '  -3.08%  '
